$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "Purpose`nid=""Par1"">To assess the management and safety of epidural or general anesthesia for Cesarean delivery in parturients with coronavirus disease (COVID-19) and their newborns, and to evaluate the standardized procedures for protecting medical staff.`n`n`nMethods`nid=""Par2"">We retrospectively reviewed the cases of parturients diagnosed with severe acute respiratory syndrome coronavirus (SARS-CoV-2) infection disease (COVID-19).`n`n Their epidemiologic history, chest computed tomography scans, laboratory measurements, and SARS-CoV-2 nucleic acid positivity were evaluated.`n`n We also recorded the patients’ demographic and clinical characteristics, anesthesia and surgery-related data, maternal and neonatal complications, as well as the health status of the involved medical staff.`n`n`nResults`nid=""Par3"">The clinical characteristics of 17 pregnant women infected with SARS-CoV-2 were similar to those previously reported in non-pregnant adult patients.`n`n All of the 17 patients underwent Cesarean delivery with anesthesia performed according to standardized anesthesia/surgery procedures.`n`n Fourteen of the patients underwent continuous epidural anesthesia with 12 experiencing significant intraoperative hypotension.`n`n Three patients received general anesthesia with tracheal intubation because emergency surgery was needed.`n`n Three of the parturients are still recovering from their Cesarean delivery and are receiving in-hospital treatment for COVID-19. Three neonates were born prematurely.`n`n There were no deaths or serious neonatal asphyxia events.`n`n All neonatal SARS-CoV-2 nucleic acid tests were negative.`n`n No medical staff were infected throughout the patient care period.`n`n`nConclusions`nid=""Par4"">Both epidural and general anesthesia were safely used for Cesarean delivery in the parturients with COVID-19. Nevertheless, the incidence of hypotension during epidural anesthesia appeared excessive.`n`n Proper patient transfer, medical staff access procedures, and effective biosafety precautions are important to protect medical staff from COVID-19.`n"
$ws.Range("E4").Value = "[Rong%Chen%NULL%0,  Yuan%Zhang%NULL%0,  Lei%Huang%NULL%0,  Bi-heng%Cheng%NULL%0,  Zhong-yuan%Xia%NULL%0,  Qing-tao%Meng%mengqingtao2018@126.com%0]"
$ws.Rows(4).AutoFit()

$ws.Range("D6").Value = "id=""Par1"">An outbreak of novel coronavirus pneumonia occurred worldwide since December 2019, which had been named COVID-19 subsequently.`n`n It is extremely transmissive that infection in pregnant women were unavoidable.`n`n The delivery process will produce large amount of contaminated media, leaving a challenge for medical personnel to ensure both the safety of the mother and infant and good self-protection.`n`n Only rare cases of pregnant women with COVID-19 are available for reference.`n`n Here, we report a 30-year-old woman had reverse transcription polymerase chain reaction-confirmed COVID-19 at 36 weeks 2 days of gestation.`n`n Significant low and high variability of fetal heart rate baseline and severe variable decelerations were repeated after admission.`n`n An emergency cesarean section at 37 weeks 1 day of gestation under combined spinal and epidural anesthesia was performed with strict protection for all personnel.`n`n Anesthesia and operation went uneventfully.`n`n None of the participants were infected.`n`n We can conclude that when confronted with cesarean section in parturient with COVID-19, careful planning and detailed preparation can improve the safety of the mother and infant and reduce the risk of infection for medical staff to help preventing and controlling the epidemic.`n`n"
$ws.Range("E6").Value = "[Yin%Du%NULL%1,  Long%Wang%NULL%1,  Gang%Wu%NULL%1,  Xiaoming%Lei%NULL%1,  Wei%Li%NULL%0,  Jianrui%Lv%lvjianrui@163.com%1]"
$ws.Rows(6).AutoFit()

$ws.Range("D8").Value = "Background`nid=""Par1"">The concerns of the highly contagious and morbid nature of Coronavirus Disease-2019 (COVID-2019) have prompted healthcare workers to implement strict droplet and contact isolation precautions.`n`n Unfortunately, some patients who may be or presumptively or confirmed as infected with COVID-2019 may also require emergent surgical procedures.`n`n As such, given the high-risk for exposure of many healthcare workers involved the complex requirements for appropriate isolation must be adhered to.`n`n`nCase presentation`nid=""Par2"">We present our experience with a 77-year-old who required emergency cardiac surgery for a presumed acute aortic syndrome in the setting of a presumed, and eventually confirmed, COVID-2019 infection.`n`n We outline the necessary steps to maintain strict isolation precautions to limit potential exposure to the surgical Team.`n`n`nConclusions`nid=""Par3"">We hereby provide our algorithm for emergent surgical procedures in critically-ill patients with presumptive or confirmed infection with COVID-2019. The insights from this case report can potentially be templated to other facilities in order to uphold high standards of infection prevention and patient safety in surgery during the current COVID-19 pandemic.`n`n`n"
$ws.Range("E8").Value = "[Michael S.%Firstenberg%msfirst@gmail.com%1,  Matthew%Libby%mnlibby@gmail.com%1,  Michael%Ochs%Michael.Ochs@usap.com%1,  Jennifer%Hanna%hannaj21@hotmail.com%1,  Julie E.%Mangino%Jemangino1@gmail.com%1,  Joseph%Forrester%Joe.Forrester@ccpcmd.com%1]"
$ws.Rows(8).AutoFit()

